# This script applies updated transition-probability values to the
# "Starting_State" matrix on the active worksheet (Portland_A).
# The workbook models game-state transition probabilities (rows sum to 1).
# Values below reflect a refreshed simulation run with more games played,
# per the accompanying commit message ("added more games, sped up
# simulate game logic, and drafted optimization logic").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1846689895470383
$ws.Cells.Item(2, 3).Value = 0.5592334494773519
$ws.Cells.Item(2, 10).Value = 0.01219512195121951
$ws.Cells.Item(2, 16).Value = 0.1393728222996516
$ws.Cells.Item(2, 19).Value = 0.1045296167247387
$ws.Cells.Item(3, 2).Value = 0.006211180124223602
$ws.Cells.Item(3, 3).Value = 0.003105590062111801
$ws.Cells.Item(3, 10).Value = 0.0124223602484472
$ws.Cells.Item(3, 16).Value = 0.7236024844720497
$ws.Cells.Item(3, 19).Value = 0.2546583850931677
$ws.Cells.Item(4, 10).Value = 0.04347826086956522
$ws.Cells.Item(4, 16).Value = 0.6956521739130435
$ws.Cells.Item(4, 19).Value = 0.2608695652173913
$ws.Cells.Item(6, 2).Value = 0.06345733041575492
$ws.Cells.Item(6, 4).Value = 0.01969365426695843
$ws.Cells.Item(6, 6).Value = 0.0700218818380744
$ws.Cells.Item(6, 10).Value = 0.2888402625820569
$ws.Cells.Item(6, 15).Value = 0.02188183807439825
$ws.Cells.Item(6, 17).Value = 0.1444201312910285
$ws.Cells.Item(6, 18).Value = 0.06345733041575492
$ws.Cells.Item(6, 19).Value = 0.3282275711159737
$ws.Cells.Item(7, 2).Value = 0.1257668711656442
$ws.Cells.Item(7, 4).Value = 0.02760736196319018
$ws.Cells.Item(7, 6).Value = 0.03680981595092025
$ws.Cells.Item(7, 10).Value = 0.1226993865030675
$ws.Cells.Item(7, 15).Value = 0.02760736196319018
$ws.Cells.Item(7, 17).Value = 0.1564417177914111
$ws.Cells.Item(7, 18).Value = 0.0705521472392638
$ws.Cells.Item(7, 19).Value = 0.4325153374233129
$ws.Cells.Item(8, 2).Value = 0.09232480533926585
$ws.Cells.Item(8, 4).Value = 0.02669632925472748
$ws.Cells.Item(8, 5).Value = 0.001112347052280311
$ws.Cells.Item(8, 6).Value = 0.06451612903225806
$ws.Cells.Item(8, 10).Value = 0.09899888765294772
$ws.Cells.Item(8, 15).Value = 0.01779755283648498
$ws.Cells.Item(8, 17).Value = 0.1635150166852058
$ws.Cells.Item(8, 18).Value = 0.1078976640711902
$ws.Cells.Item(8, 19).Value = 0.4271412680756396
$ws.Cells.Item(9, 2).Value = 0.07276507276507277
$ws.Cells.Item(9, 4).Value = 0.01871101871101871
$ws.Cells.Item(9, 6).Value = 0.05405405405405406
$ws.Cells.Item(9, 10).Value = 0.1164241164241164
$ws.Cells.Item(9, 15).Value = 0.02702702702702703
$ws.Cells.Item(9, 17).Value = 0.1995841995841996
$ws.Cells.Item(9, 18).Value = 0.103950103950104
$ws.Cells.Item(9, 19).Value = 0.4074844074844075
$ws.Cells.Item(10, 2).Value = 0.1081081081081081
$ws.Cells.Item(10, 4).Value = 0.01954261954261954
$ws.Cells.Item(10, 5).Value = 0.0004158004158004158
$ws.Cells.Item(10, 6).Value = 0.06278586278586279
$ws.Cells.Item(10, 10).Value = 0.1239085239085239
$ws.Cells.Item(10, 15).Value = 0.01455301455301455
$ws.Cells.Item(10, 17).Value = 0.1808731808731809
$ws.Cells.Item(10, 18).Value = 0.09022869022869023
$ws.Cells.Item(10, 19).Value = 0.3995841995841996
$ws.Cells.Item(11, 7).Value = 0.150390625
$ws.Cells.Item(11, 10).Value = 0.109375
$ws.Cells.Item(11, 11).Value = 0.208984375
$ws.Cells.Item(11, 12).Value = 0.5078125
$ws.Cells.Item(11, 19).Value = 0.0234375
$ws.Cells.Item(12, 7).Value = 0.7323420074349443
$ws.Cells.Item(12, 10).Value = 0.1821561338289963
$ws.Cells.Item(12, 11).Value = 0.01486988847583643
$ws.Cells.Item(12, 12).Value = 0.02973977695167286
$ws.Cells.Item(12, 19).Value = 0.04089219330855019
$ws.Cells.Item(13, 6).Value = 0.01136363636363636
$ws.Cells.Item(13, 7).Value = 0.6363636363636364
$ws.Cells.Item(13, 10).Value = 0.2727272727272727
$ws.Cells.Item(13, 19).Value = 0.07954545454545454
$ws.Cells.Item(14, 7).Value = 0.6
$ws.Cells.Item(14, 10).Value = 0.4
$ws.Cells.Item(15, 6).Value = 0.03267973856209151
$ws.Cells.Item(15, 8).Value = 0.1546840958605664
$ws.Cells.Item(15, 9).Value = 0.07625272331154684
$ws.Cells.Item(15, 10).Value = 0.3464052287581699
$ws.Cells.Item(15, 11).Value = 0.0457516339869281
$ws.Cells.Item(15, 13).Value = 0.01525054466230937
$ws.Cells.Item(15, 15).Value = 0.07407407407407407
$ws.Cells.Item(15, 19).Value = 0.2549019607843137
$ws.Cells.Item(16, 6).Value = 0.0273972602739726
$ws.Cells.Item(16, 8).Value = 0.1753424657534247
$ws.Cells.Item(16, 9).Value = 0.0821917808219178
$ws.Cells.Item(16, 10).Value = 0.3972602739726027
$ws.Cells.Item(16, 11).Value = 0.09315068493150686
$ws.Cells.Item(16, 13).Value = 0.02465753424657534
$ws.Cells.Item(16, 14).Value = 0.005479452054794521
$ws.Cells.Item(16, 15).Value = 0.06027397260273973
$ws.Cells.Item(16, 19).Value = 0.1342465753424658
$ws.Cells.Item(17, 6).Value = 0.02010050251256281
$ws.Cells.Item(17, 8).Value = 0.1947236180904523
$ws.Cells.Item(17, 9).Value = 0.10678391959799
$ws.Cells.Item(17, 10).Value = 0.3881909547738693
$ws.Cells.Item(17, 11).Value = 0.08668341708542714
$ws.Cells.Item(17, 13).Value = 0.01633165829145729
$ws.Cells.Item(17, 14).Value = 0.001256281407035176
$ws.Cells.Item(17, 15).Value = 0.06407035175879397
$ws.Cells.Item(17, 19).Value = 0.1218592964824121
$ws.Cells.Item(18, 6).Value = 0.04830917874396135
$ws.Cells.Item(18, 8).Value = 0.1714975845410628
$ws.Cells.Item(18, 9).Value = 0.1328502415458937
$ws.Cells.Item(18, 10).Value = 0.3236714975845411
$ws.Cells.Item(18, 11).Value = 0.0966183574879227
$ws.Cells.Item(18, 13).Value = 0.00966183574879227
$ws.Cells.Item(18, 14).Value = 0.002415458937198068
$ws.Cells.Item(18, 15).Value = 0.08695652173913043
$ws.Cells.Item(18, 19).Value = 0.1280193236714976
$ws.Cells.Item(19, 6).Value = 0.02453047144499809
$ws.Cells.Item(19, 8).Value = 0.2108087389804523
$ws.Cells.Item(19, 9).Value = 0.1061709467228823
$ws.Cells.Item(19, 10).Value = 0.3583748562667689
$ws.Cells.Item(19, 11).Value = 0.08738980452280568
$ws.Cells.Item(19, 13).Value = 0.02069758528171713
$ws.Cells.Item(19, 14).Value = 0.0007665772326561902
$ws.Cells.Item(19, 15).Value = 0.0697585281717133
$ws.Cells.Item(19, 19).Value = 0.1215024913760061
